$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily balance rows (9-13), mirroring the existing pattern of
# Date / ValueLTC / IncrementDaily (calculated table column formula).
$ws.Range("A9").Value = 44199
$ws.Range("B9").Value = 0.00149578
$ws.Range("C9").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B8,0),8)"

$ws.Range("A10").Value = 44200
$ws.Range("B10").Value = 0.00162743
$ws.Range("C10").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B9,0),8)"

$ws.Range("A11").Value = 44201
$ws.Range("B11").Value = 0.0017491
$ws.Range("C11").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B10,0),8)"

$ws.Range("A12").Value = 44202
$ws.Range("B12").Value = 0.00187205
$ws.Range("C12").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B11,0),8)"

$ws.Range("A13").Value = 44203
$ws.Range("B13").Value = 0.00197919
$ws.Range("C13").Formula = "=ROUND(IFERROR(BalanceDaily[[#This Row],[ValueLTC]]-B12,0),8)"

# Extend the BalanceDaily table to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C13"))

# Reflect the last interaction: selecting the whole of column D
# (as if preparing/checking a condition before the next log entry).
$ws.Columns("D").Select()
